$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 5435.364
$ws.Range("J17").Value = 5670.381
$ws.Range("L17").Value = 17011.143
$ws.Range("N17").Value = -17347.143

$ws.Range("H40").Value = 2530.05
$ws.Range("I40").Value = 3570.4285
$ws.Range("J40").Value = 1969.8462
$ws.Range("K40").Value = 3570.4285
$ws.Range("L40").Value = 1969.8462
$ws.Range("M40").Value = -3395.4285
$ws.Range("N40").Value = -2319.8462

$ws.Range("H129").Value = 1221.7778
$ws.Range("I129").Value = 498.5
$ws.Range("J129").Value = 1428.4286
$ws.Range("K129").Value = 1495.5
$ws.Range("L129").Value = 4285.2858
$ws.Range("M129").Value = 3504.5
$ws.Range("N129").Value = -14285.2858

$ws.Range("H137").Value = 2473.6064
$ws.Range("I137").Value = 1340.8889
$ws.Range("K137").Value = 4022.6667
$ws.Range("M137").Value = -1472.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7001
$ws.Range("I32").Value = 5113.5
$ws.Range("J32").Value = 19012.363
$ws.Range("K32").Value = 5113.5
$ws.Range("L32").Value = 19012.363
$ws.Range("M32").Value = -4826.5
$ws.Range("N32").Value = -19586.363

$ws.Range("H61").Value = 4287.077
$ws.Range("I61").Value = 3122.139
$ws.Range("K61").Value = 3122.139
$ws.Range("M61").Value = -2910.139

$ws.Range("H132").Value = 6292.6113
$ws.Range("I132").Value = 1983.9166
$ws.Range("J132").Value = 8446.958000000001
$ws.Range("K132").Value = 5951.7498
$ws.Range("L132").Value = 25340.874
$ws.Range("M132").Value = -3421.7498
$ws.Range("N132").Value = -30400.874

$ws.Range("H135").Value = 53937.785
$ws.Range("J135").Value = 53937.785
$ws.Range("L135").Value = 53937.785
$ws.Range("N135").Value = -64077.785

$ws.Range("H136").Value = 4287.077
$ws.Range("I136").Value = 3122.139
$ws.Range("K136").Value = 9366.417000000001
$ws.Range("M136").Value = -6816.417000000001

$ws.Range("H139").Value = 52000
$ws.Range("J139").Value = 52000
$ws.Range("L139").Value = 52000
$ws.Range("N139").Value = -62280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2836.6428
$ws.Range("I107").Value = 2750
$ws.Range("J107").Value = 3053.25
$ws.Range("K107").Value = 2750
$ws.Range("L107").Value = 3053.25
$ws.Range("M107").Value = -830
$ws.Range("N107").Value = -6893.25

$ws.Range("H134").Value = 6400.6523
$ws.Range("I134").Value = 6011.8
$ws.Range("J134").Value = 8993
$ws.Range("K134").Value = 18035.4
$ws.Range("L134").Value = 26979
$ws.Range("M134").Value = -15500.4
$ws.Range("N134").Value = -32049

$ws.Range("H141").Value = 43596.668
$ws.Range("J141").Value = 43596.668
$ws.Range("L141").Value = 43596.668
$ws.Range("N141").Value = -53956.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2739.283
$ws.Range("I31").Value = 2394.8
$ws.Range("J31").Value = 3046.8572
$ws.Range("K31").Value = 2394.8
$ws.Range("L31").Value = 3046.8572
$ws.Range("M31").Value = -2099.8
$ws.Range("N31").Value = -3636.8572

$ws.Range("H34").Value = 2739.283
$ws.Range("I34").Value = 2394.8
$ws.Range("J34").Value = 3046.8572
$ws.Range("K34").Value = 2394.8
$ws.Range("L34").Value = 3046.8572
$ws.Range("M34").Value = -2192.8
$ws.Range("N34").Value = -3450.8572

$ws.Range("H107").Value = 876.8461
$ws.Range("I107").Value = 1194.4286
$ws.Range("J107").Value = 506.33334
$ws.Range("K107").Value = 1194.4286
$ws.Range("L107").Value = 506.33334
$ws.Range("M107").Value = 725.5714
$ws.Range("N107").Value = -4346.33334

$ws.Range("H112").Value = 67767.336
$ws.Range("J112").Value = 67767.336
$ws.Range("L112").Value = 67767.336
$ws.Range("N112").Value = -70721.336

$ws.Range("H122").Value = 21169.555
$ws.Range("I122").Value = 9653
$ws.Range("J122").Value = 30382.8
$ws.Range("K122").Value = 28959
$ws.Range("L122").Value = 91148.39999999999
$ws.Range("M122").Value = -26509
$ws.Range("N122").Value = -96048.39999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 13892404
$ws.Range("I5").Value = 522.4231
$ws.Range("J5").Value = 50011296
$ws.Range("K5").Value = 1567.2693
$ws.Range("L5").Value = 150033888
$ws.Range("M5").Value = -1455.2693
$ws.Range("N5").Value = -150034112

$ws.Range("H131").Value = 48981.75
$ws.Range("I131").Value = 2061.6667
$ws.Range("J131").Value = 87370.91
$ws.Range("K131").Value = 6185.000100000001
$ws.Range("L131").Value = 262112.73
$ws.Range("M131").Value = -1145.000100000001
$ws.Range("N131").Value = -272192.73

$ws.Range("H135").Value = 13892404
$ws.Range("I135").Value = 522.4231
$ws.Range("J135").Value = 50011296
$ws.Range("K135").Value = 4701.8079
$ws.Range("L135").Value = 450101664
$ws.Range("M135").Value = -2166.8079
$ws.Range("N135").Value = -450106734

$ws.Range("H137").Value = 31285320
$ws.Range("I137").Value = 50001428
$ws.Range("J137").Value = 91805.336
$ws.Range("K137").Value = 150004284
$ws.Range("L137").Value = 275416.008
$ws.Range("M137").Value = -149999184
$ws.Range("N137").Value = -285616.008

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2264.4167
$ws.Range("I132").Value = 1540.5
$ws.Range("J132").Value = 2781.5
$ws.Range("K132").Value = 4621.5
$ws.Range("L132").Value = 8344.5
$ws.Range("M132").Value = -2091.5
$ws.Range("N132").Value = -13404.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 791.4167
$ws.Range("I46").Value = 750
$ws.Range("J46").Value = 799.7
$ws.Range("K46").Value = 750
$ws.Range("L46").Value = 799.7
$ws.Range("M46").Value = -562
$ws.Range("N46").Value = -1175.7

$ws.Range("H122").Value = 5605.645
$ws.Range("I122").Value = 4642.609
$ws.Range("J122").Value = 8374.375
$ws.Range("K122").Value = 13927.827
$ws.Range("L122").Value = 25123.125
$ws.Range("M122").Value = -11477.827
$ws.Range("N122").Value = -30023.125

$ws.Range("H136").Value = 6709.448
$ws.Range("I136").Value = 5055.2856
$ws.Range("J136").Value = 8253.333000000001
$ws.Range("K136").Value = 15165.8568
$ws.Range("L136").Value = 24759.999
$ws.Range("M136").Value = -12615.8568
$ws.Range("N136").Value = -29859.999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 38249
$ws.Range("J63").Value = 38249
$ws.Range("L63").Value = 38249
$ws.Range("N63").Value = -39497

$ws.Range("H66").Value = 38249
$ws.Range("J66").Value = 38249
$ws.Range("L66").Value = 114747
$ws.Range("N66").Value = -120987

$ws.Range("H122").Value = 2997.182
$ws.Range("I122").Value = 2633
$ws.Range("J122").Value = 3968.3333
$ws.Range("K122").Value = 7899
$ws.Range("L122").Value = 11904.9999
$ws.Range("M122").Value = -5449
$ws.Range("N122").Value = -16804.9999

$ws.Range("H132").Value = 2129.8333
$ws.Range("I132").Value = 1983.3529
$ws.Range("J132").Value = 2485.5715
$ws.Range("K132").Value = 5950.0587
$ws.Range("L132").Value = 7456.7145
$ws.Range("M132").Value = -3420.0587
$ws.Range("N132").Value = -12516.7145
